# The post "「化学反応」" (row 189) was removed from the blog export.
# Delete the entire row 189; all following rows (190-216) shift up by one,
# and the sheet dimension shrinks from A1:C216 to A1:C215 automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(189).Delete()
